# Generate Report for Handback
# Updates the localization-status workbook to reflect the handback that just
# completed: status text, handback datetimes, and the newly-produced target /
# handback files (with hyperlinks) for both the zh-cn and de-de locales.

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

$target291882 = "291882cc-0c61-4895-b9d5-ef6e7719649e.md"
$target291882Url = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3a0dc0fc0a30bafe838398a999fbb0467b2ab0c0/e2e/291882cc-0c61-4895-b9d5-ef6e7719649e.md"
$targetEebc27 = "eebc27f0-c2c4-438b-90ba-0dccc57d4328.md"
$targetEebc27Url = "https://github.com/OpenLocalizationTestOrg/ol-test4/blob/3a0dc0fc0a30bafe838398a999fbb0467b2ab0c0/e2e/eebc27f0-c2c4-438b-90ba-0dccc57d4328.md"

$zhHandback291882 = "291882cc-0c61-4895-b9d5-ef6e7719649e.d4a4d18925de74e1bcb6644e994b090d74de43bc.zh-cn.xlf"
$zhHandbackEebc27 = "eebc27f0-c2c4-438b-90ba-0dccc57d4328.ad4b3f1333a07192ef83422400b8429e4c63c05b.zh-cn.xlf"
$deHandback291882 = "291882cc-0c61-4895-b9d5-ef6e7719649e.d4a4d18925de74e1bcb6644e994b090d74de43bc.de-de.xlf"
$deHandbackEebc27 = "eebc27f0-c2c4-438b-90ba-0dccc57d4328.ad4b3f1333a07192ef83422400b8429e4c63c05b.de-de.xlf"

$zhHandbackTime = "2017-02-22 08:35:47"
$deHandbackTime = "2017-02-22 08:36:10"

# Hyperlink font look-and-feel (matches the existing custom "HyperLink" style
# already used on column A: underlined, #6495ED).
$hyperlinkColor = 15570276  # BGR long for RGB 0x6495ED
function Set-HyperlinkLook($range) {
    $range.Font.Underline = 2
    $range.Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------------
# Overview sheet — widen the per-locale status columns and refresh the status
# text shown for both files in both locale columns.
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value2 = $statusText
$wsOverview.Range("F2").Value2 = $statusText
$wsOverview.Range("E3").Value2 = $statusText
$wsOverview.Range("F3").Value2 = $statusText
$wsOverview.Columns.Item(5).ColumnWidth = 29.15
$wsOverview.Columns.Item(6).ColumnWidth = 29.15

# ---------------------------------------------------------------------------
# zh-cn sheet — record the generated target file + handback xliff for each
# row, stamp the handback datetime, refresh status text, and widen columns.
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")

$wsZh.Range("C2").Value2 = $statusText
$wsZh.Range("C3").Value2 = $statusText

$wsZh.Range("J2").Value2 = $target291882
$wsZh.Range("K2").Value2 = $zhHandback291882
$wsZh.Range("L2").Value2 = $zhHandbackTime

$wsZh.Range("J3").Value2 = $targetEebc27
$wsZh.Range("K3").Value2 = $zhHandbackEebc27
$wsZh.Range("L3").Value2 = $zhHandbackTime

$wsZh.Hyperlinks.Add($wsZh.Range("J2"), $target291882Url, $null, $null, $target291882) | Out-Null
$wsZh.Hyperlinks.Add($wsZh.Range("J3"), $targetEebc27Url, $null, $null, $targetEebc27) | Out-Null
Set-HyperlinkLook $wsZh.Range("J2")
Set-HyperlinkLook $wsZh.Range("J3")

$wsZh.Columns.Item(3).ColumnWidth = 29.15
$wsZh.Columns.Item(10).ColumnWidth = 39.15
$wsZh.Columns.Item(11).ColumnWidth = 39.15

# ---------------------------------------------------------------------------
# de-de sheet — same shape of update, using the de-de handoff/handback file
# names and the de-de handback timestamp.
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")

$wsDe.Range("C2").Value2 = $statusText
$wsDe.Range("C3").Value2 = $statusText

$wsDe.Range("J2").Value2 = $target291882
$wsDe.Range("K2").Value2 = $deHandback291882
$wsDe.Range("L2").Value2 = $deHandbackTime

$wsDe.Range("J3").Value2 = $targetEebc27
$wsDe.Range("K3").Value2 = $deHandbackEebc27
$wsDe.Range("L3").Value2 = $deHandbackTime

$wsDe.Hyperlinks.Add($wsDe.Range("J2"), $target291882Url, $null, $null, $target291882) | Out-Null
$wsDe.Hyperlinks.Add($wsDe.Range("J3"), $targetEebc27Url, $null, $null, $targetEebc27) | Out-Null
Set-HyperlinkLook $wsDe.Range("J2")
Set-HyperlinkLook $wsDe.Range("J3")

$wsDe.Columns.Item(3).ColumnWidth = 29.15
$wsDe.Columns.Item(10).ColumnWidth = 39.15
$wsDe.Columns.Item(11).ColumnWidth = 39.15
